# Applies the cryptos.xlsx price/volume update described by the commit
# "Updated cryptos list on Wed Aug 23 20:55:54 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is numeric-looking (e.g. "1.002", "217.71") must be
# force-formatted as Text first, otherwise Excel auto-converts the assigned
# string into a real number on Range.Value assignment (losing exact text,
# e.g. trailing zeros such as "0.8650" -> 0.865).
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D10",
    "D11",
    "D14",
    "D16",
    "D18",
    "D19",
    "D20",
    "D24",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D35",
    "D36",
    "D39",
    "D41",
    "D43",
    "D47",
    "D49",
    "D50",
    "D51",
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new cell values (row order follows the sheet, columns B-E).
$ws.Range("D2").Value = "26.635.28"
$ws.Range("E2").Value = "  +2.87%  "
$ws.Range("D3").Value = "1.688.00"
$ws.Range("E3").Value = "  +3.46%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "217.71"
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").Value = "0.5341"
$ws.Range("E6").Value = "  +2.72%  "
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "0.2683"
$ws.Range("E8").Value = "  +4.57%  "
$ws.Range("E9").Value = "  +3.41%  "
$ws.Range("D10").Value = "21.69"
$ws.Range("E10").Value = "  +7.10%  "
$ws.Range("D11").Value = "0.07793"
$ws.Range("E11").Value = "  +3.12%  "
$ws.Range("D12").Value = "1.689.69"
$ws.Range("E12").Value = "  +3.57%  "
$ws.Range("E13").Value = "  +3.60%  "
$ws.Range("D14").Value = "0.5625"
$ws.Range("E14").Value = "  +3.97%  "
$ws.Range("D15").Value = "0.0₅8448"
$ws.Range("E15").Value = "  +6.94%  "
$ws.Range("D16").Value = "66.27"
$ws.Range("E16").Value = "  +2.94%  "
$ws.Range("D17").Value = "26.684.19"
$ws.Range("E17").Value = "  +3.04%  "
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "4.813"
$ws.Range("E19").Value = "  +4.60%  "
$ws.Range("D20").Value = "195.60"
$ws.Range("E20").Value = "  +6.22%  "
$ws.Range("E21").Value = "  +4.29%  "
$ws.Range("E22").Value = "  +5.48%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "144.14"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "0.1292"
$ws.Range("E25").Value = "  +7.58%  "
$ws.Range("D26").Value = "7.491"
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("E27").Value = "  +5.34%  "
$ws.Range("D28").Value = "1.429"
$ws.Range("E28").Value = "  +4.42%  "
$ws.Range("D29").Value = "0.06172"
$ws.Range("E29").Value = "  +4.13%  "
$ws.Range("D30").Value = "1.282"
$ws.Range("D31").Value = "3.606"
$ws.Range("E31").Value = "  +7.78%  "
$ws.Range("D32").Value = "3.476"
$ws.Range("E32").Value = "  +3.97%  "
$ws.Range("E33").Value = "  +6.39%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.427"
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "2.798"
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("E38").Value = "  +3.78%  "
$ws.Range("D39").Value = "6.026"
$ws.Range("E39").Value = "  +6.79%  "
$ws.Range("D40").Value = "1.080.53"
$ws.Range("E40").Value = "  +5.97%  "
$ws.Range("D41").Value = "0.8650"
$ws.Range("E41").Value = "  +3.30%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "100.52"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D44").Value = "1.839.45"
$ws.Range("E44").Value = "  +3.15%  "
$ws.Range("D45").Value = "0.0₈109"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("E46").Value = "  +6.13%  "
$ws.Range("D47").Value = "8.218"
$ws.Range("E47").Value = "  +3.65%  "
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").Value = "0.05225"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").Value = "6.111"
$ws.Range("E50").Value = "  +5.87%  "
$ws.Range("D51").Value = "0.4243"
$ws.Range("E51").Value = "  +0.34%  "
